$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value. All values are written as literal
# text (matching the source inlineStr cells) by temporarily forcing a Text
# number format, assigning the value, then clearing the format back to General
# so numeric-looking strings (e.g. "1.003", "116.30", "0.000008657") are not
# auto-coerced into numbers/scientific notation by Excel.
$updates = [ordered]@{
    'D2' = '27.148.03'
    'E2' = '  -0.15%  '
    'D3' = '1.899.45'
    'E3' = '  -0.17%  '
    'E4' = '  +0.25%  '
    'D5' = '307.24'
    'E5' = '  +0.18%  '
    'E6' = '  +0.26%  '
    'D7' = '0.5238'
    'E7' = '  +0.07%  '
    'D8' = '0.3805'
    'E8' = '  +0.98%  '
    'D9' = '0.07304'
    'E9' = '  +0.86%  '
    'D10' = '21.35'
    'E10' = '  +0.81%  '
    'D11' = '0.9064'
    'E11' = '  +0.88%  '
    'D12' = '0.08199'
    'E12' = '  -2.61%  '
    'D13' = '1.879.60'
    'E13' = '  -1.18%  '
    'D14' = '95.23'
    'E14' = '  +0.70%  '
    'D15' = '5.344'
    'E15' = '  +1.45%  '
    'D16' = '1.003'
    'E16' = '  +0.17%  '
    'D17' = '0.000008657'
    'E17' = '  +0.80%  '
    'D18' = '14.67'
    'E18' = '  +1.28%  '
    'D20' = '27.194.72'
    'E20' = '  -0.13%  '
    'D21' = '5.115'
    'E21' = '  +1.21%  '
    'D22' = '2.123.06'
    'E22' = '  -0.89%  '
    'D23' = '10.78'
    'E23' = '  +1.89%  '
    'D24' = '6.476'
    'E24' = '  +0.94%  '
    'D25' = '2.338'
    'E25' = '  +2.28%  '
    'D26' = '149.59'
    'E26' = '  +2.01%  '
    'D27' = '18.27'
    'E27' = '  +0.72%  '
    'D28' = '1.743'
    'E28' = '  -0.55%  '
    'D29' = '115.41'
    'E29' = '  +0.52%  '
    'D30' = '4.826'
    'E30' = '  +0.81%  '
    'D31' = '4.856'
    'E31' = '  -1.39%  '
    'D32' = '0.09246'
    'E32' = '  +0.38%  '
    'D33' = '0.05052'
    'E33' = '  -0.08%  '
    'D34' = '0.7944'
    'E34' = '  -2.40%  '
    'D35' = '1.222'
    'E35' = '  -1.34%  '
    'D36' = '2.958'
    'E36' = '  +0.30%  '
    'D37' = '3.388'
    'E37' = '  +0.30%  '
    'D38' = '2.672'
    'E38' = '  +4.25%  '
    'D39' = '0.5732'
    'E39' = '  +0.70%  '
    'E40' = '  +0.95%  '
    'D41' = '1.079'
    'E41' = '  +0.64%  '
    'D42' = '9.011'
    'D43' = '6.614'
    'E43' = '  -0.39%  '
    'D44' = '116.30'
    'E44' = '  -1.66%  '
    'D45' = '0.1518'
    'E45' = '  +0.43%  '
    'D46' = '0.4895'
    'E46' = '  +1.53%  '
    'E47' = '  +0.23%  '
    'D48' = '10.18'
    'E48' = '  +0.29%  '
    'D49' = '1.638'
    'E49' = '  +1.59%  '
    'D50' = '38.56'
    'E50' = '  +3.09%  '
    'D51' = '64.08'
    'E51' = '  +0.77%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.ClearFormats()
}

Write-Host "Updated" $updates.Count "cells"
